$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The three data rows (2, 3, 4) get their values cyclically rotated:
#   new row2 = old row4 values
#   new row3 = old row2 values
#   new row4 = old row3 values
# Column R (Origen) for row4 stays "Provincia de Limarí" both before and after,
# and column Q (Unidad de comercialización) is identical across all rows, so no
# change is required there - but we set every touched field explicitly for clarity
# and to exactly match the target values.

# Row 2 (becomes old row 4's data)
$ws.Range("D2").Value = 44320
$ws.Range("M2").Value = 50
$ws.Range("N2").Value = 18000
$ws.Range("O2").Value = 20000
$ws.Range("P2").Value = 18800
$ws.Range("R2").Value = "Provincia de Limarí"
$ws.Range("S2").Value = 1044

# Row 3 (becomes old row 2's data)
$ws.Range("D3").Value = 44362
$ws.Range("M3").Value = 100
$ws.Range("N3").Value = 19000
$ws.Range("O3").Value = 20000
$ws.Range("P3").Value = 19500
$ws.Range("R3").Value = "Provincia de Curicó"
$ws.Range("S3").Value = 1083

# Row 4 (becomes old row 3's data)
$ws.Range("D4").Value = 44719
$ws.Range("M4").Value = 50
$ws.Range("N4").Value = 20000
$ws.Range("O4").Value = 21000
$ws.Range("P4").Value = 20400
$ws.Range("R4").Value = "Provincia de Limarí"
$ws.Range("S4").Value = 1133
